# Fruta / hortaliza, semanal
# The weekly refresh re-shuffles the daily price rows (2..47): each row keeps
# its constant columns (Mercado, Region, Codreg, Categoria, Calidad, Unidad,
# Kg o Unidades, Clasificacion) but receives a "new" date / volume / price /
# origin combination, effectively a permutation of the existing rows.
# Build that permutation explicitly (targetRow -> sourceRow, values taken
# from the ORIGINAL, pre-edit sheet) and apply it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowMap = @{
    2 = 28
    3 = 24
    4 = 32
    5 = 43
    6 = 2
    7 = 26
    8 = 14
    9 = 35
    10 = 12
    11 = 38
    12 = 13
    13 = 18
    14 = 47
    15 = 42
    16 = 6
    17 = 30
    18 = 23
    19 = 44
    20 = 11
    21 = 20
    22 = 7
    23 = 41
    24 = 16
    25 = 22
    26 = 31
    27 = 45
    28 = 33
    29 = 4
    30 = 19
    31 = 5
    32 = 37
    33 = 29
    34 = 27
    35 = 39
    36 = 8
    37 = 9
    38 = 40
    39 = 15
    40 = 25
    41 = 21
    42 = 36
    43 = 3
    44 = 46
    45 = 34
    46 = 17
    47 = 10
}

# Columns that actually vary between rows and therefore move with the
# permutation: D(4) Fecha, H(8) Variedad, J(10) Volumen, K(11) Precio minimo,
# L(12) Precio maximo, M(13) Precio promedio ponderado, O(15) Origen,
# P(16) Precio $/Kg.
$cols = @(4, 8, 10, 11, 12, 13, 15, 16)

# Snapshot every needed source cell BEFORE any writes, since sources and
# targets overlap (it's a permutation).
$snapshot = @{}
foreach ($r in $rowMap.Keys) {
    $src = $rowMap[$r]
    if (-not $snapshot.ContainsKey($src)) {
        $vals = @{}
        foreach ($c in $cols) {
            $vals[$c] = $ws.Cells.Item($src, $c).Value()
        }
        $snapshot[$src] = $vals
    }
}

foreach ($r in $rowMap.Keys) {
    $src = $rowMap[$r]
    $vals = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $vals[$c]
    }
}
